$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "0.844 (0.004)"
$ws.Range("C3").Value = "0.877 (0.007)"
$ws.Range("C4").Value = "0.932 (0.007)"
$ws.Range("C5").Value = "0.953 (0.003)"
$ws.Range("C6").Value = "0.999 (0.001)"
$ws.Range("C7").Value = "0.999 (0.001)"
$ws.Range("C8").Value = "1.000 (0.000)"
$ws.Range("C9").Value = "1.000 (0.000)"
